$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update left table (A:H), rows 3-6 with new crude/crisis/panic/sc stats ---
$ws.Range("A3").Value = "crude"
$ws.Range("B3").Value = 0.7941176470588235
$ws.Range("C3").Value = 27
$ws.Range("D3").Value = 27
$ws.Range("H3").Value = 7

$ws.Range("A4").Value = "crisis"
$ws.Range("B4").Value = 0.6164383561643836
$ws.Range("C4").Value = 180
$ws.Range("D4").Value = 180
$ws.Range("H4").Value = 112

$ws.Range("A5").Value = "panic"
$ws.Range("B5").Value = 0.1802325581395349
$ws.Range("C5").Value = 93
$ws.Range("D5").Value = 93
$ws.Range("H5").Value = 423

$ws.Range("A6").Value = "sc"
$ws.Range("B6").Value = 0.1587301587301587
$ws.Range("C6").Value = 30
$ws.Range("D6").Value = 30
$ws.Range("H6").Value = 159

# --- Row 7 of the left table ("low") no longer exists; clear it entirely ---
$ws.Range("A7:H7").Clear()

# --- Update right table (J:Q), rows 3-27 with the rerun stats (larger dataset) ---
$ws.Range("J3").Value = "best"
$ws.Range("K3").Value = 0.9491525423728814
$ws.Range("L3").Value = 56
$ws.Range("M3").Value = 56
$ws.Range("Q3").Value = 3

$ws.Range("J4").Value = "interesting"
$ws.Range("K4").Value = 0.9393939393939394
$ws.Range("L4").Value = 31
$ws.Range("M4").Value = 31
$ws.Range("Q4").Value = 2

$ws.Range("J5").Value = "love"
$ws.Range("K5").Value = 0.8913043478260869
$ws.Range("L5").Value = 41
$ws.Range("M5").Value = 41
$ws.Range("Q5").Value = 5

$ws.Range("J6").Value = "great"
$ws.Range("K6").Value = 0.8660714285714286
$ws.Range("L6").Value = 97
$ws.Range("M6").Value = 97
$ws.Range("Q6").Value = 15

$ws.Range("J7").Value = "special"
$ws.Range("K7").Value = 0.8333333333333334
$ws.Range("L7").Value = 30
$ws.Range("M7").Value = 30
$ws.Range("Q7").Value = 6

$ws.Range("J8").Value = "thank"
$ws.Range("K8").Value = 0.828125
$ws.Range("L8").Value = 106
$ws.Range("M8").Value = 106
$ws.Range("Q8").Value = 22

$ws.Range("J9").Value = "thanks"
$ws.Range("K9").Value = 0.8170731707317073
$ws.Range("L9").Value = 67
$ws.Range("M9").Value = 67
$ws.Range("Q9").Value = 15

$ws.Range("J10").Value = "positive"
$ws.Range("K10").Value = 0.7758620689655172
$ws.Range("L10").Value = 45
$ws.Range("M10").Value = 45
$ws.Range("Q10").Value = 13

$ws.Range("J11").Value = "free"
$ws.Range("K11").Value = 0.7666666666666667
$ws.Range("L11").Value = 92
$ws.Range("M11").Value = 92
$ws.Range("Q11").Value = 28

$ws.Range("J12").Value = "safety"
$ws.Range("K12").Value = 0.7254901960784313
$ws.Range("L12").Value = 37
$ws.Range("M12").Value = 37
$ws.Range("Q12").Value = 14

$ws.Range("J13").Value = "safe"
$ws.Range("K13").Value = 0.7253521126760564
$ws.Range("L13").Value = 103
$ws.Range("M13").Value = 103
$ws.Range("Q13").Value = 39

$ws.Range("J14").Value = "support"
$ws.Range("K14").Value = 0.7169811320754716
$ws.Range("L14").Value = 76
$ws.Range("M14").Value = 76
$ws.Range("Q14").Value = 30

$ws.Range("J15").Value = "good"
$ws.Range("K15").Value = 0.69375
$ws.Range("L15").Value = 111
$ws.Range("M15").Value = 111
$ws.Range("Q15").Value = 49

$ws.Range("J16").Value = "fresh"
$ws.Range("K16").Value = 0.6666666666666666
$ws.Range("L16").Value = 32
$ws.Range("M16").Value = 32
$ws.Range("Q16").Value = 16

$ws.Range("J17").Value = "better"
$ws.Range("K17").Value = 0.6349206349206349
$ws.Range("L17").Value = 40
$ws.Range("M17").Value = 40
$ws.Range("Q17").Value = 23

$ws.Range("J18").Value = "well"
$ws.Range("K18").Value = 0.6170212765957447
$ws.Range("L18").Value = 58
$ws.Range("M18").Value = 58
$ws.Range("Q18").Value = 36

$ws.Range("J19").Value = "relief"
$ws.Range("K19").Value = 0.6
$ws.Range("L19").Value = 30
$ws.Range("M19").Value = 30
$ws.Range("Q19").Value = 20

$ws.Range("J20").Value = "heroes"
$ws.Range("K20").Value = 0.5957446808510638
$ws.Range("L20").Value = 28
$ws.Range("M20").Value = 28
$ws.Range("Q20").Value = 19

$ws.Range("J21").Value = "hand"
$ws.Range("K21").Value = 0.4986945169712794
$ws.Range("L21").Value = 191
$ws.Range("M21").Value = 191
$ws.Range("Q21").Value = 192

$ws.Range("J22").Value = "like"
$ws.Range("K22").Value = 0.4705882352941176
$ws.Range("L22").Value = 160
$ws.Range("M22").Value = 160
$ws.Range("Q22").Value = 180

$ws.Range("J23").Value = "care"
$ws.Range("K23").Value = 0.4269662921348314
$ws.Range("L23").Value = 38
$ws.Range("M23").Value = 38
$ws.Range("Q23").Value = 51

$ws.Range("J24").Value = "help"
$ws.Range("K24").Value = 0.4203389830508474
$ws.Range("L24").Value = 124
$ws.Range("M24").Value = 124
$ws.Range("Q24").Value = 171

$ws.Range("J25").Value = "hope"
$ws.Range("K25").Value = 0.4
$ws.Range("L25").Value = 26
$ws.Range("M25").Value = 26
$ws.Range("Q25").Value = 39

$ws.Range("J26").Value = "protect"
$ws.Range("K26").Value = 0.3698630136986301
$ws.Range("L26").Value = 27
$ws.Range("M26").Value = 27
$ws.Range("Q26").Value = 46

$ws.Range("J27").Value = "please"
$ws.Range("K27").Value = 0.3514644351464435
$ws.Range("L27").Value = 84
$ws.Range("M27").Value = 84
$ws.Range("Q27").Value = 155
